$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update H column formulas to use PI() instead of the literal 3.1415
$ws.Range("H2").Formula = "=F2/((B2/20)^2*PI()*D2/10)"
$ws.Range("H3").Formula = "=F3/((B3/20)^2*PI()*D3/10)"
$ws.Range("H4").Formula = "=F4/((B4/20)^2*PI()*D4/10)"
$ws.Range("H7").Formula = "=F7/((B7/20)^2*PI()*D7/10)"
$ws.Range("H8").Formula = "=F8/((B8/20)^2*PI()*D8/10)"

# Add new I column formulas (Uncertainty in Density)
$ws.Range("I2").Formula = "=SQRT((C2/B2)^2+(E2/D2)^2+(G2/F2)^2)*H2"
$ws.Range("I3").Formula = "=SQRT((C3/B3)^2+(E3/D3)^2+(G3/F3)^2)*H3"
$ws.Range("I4").Formula = "=SQRT((C4/B4)^2+(E4/D4)^2+(G4/F4)^2)*H4"
$ws.Range("I7").Formula = "=SQRT((C7/B7)^2+(E7/D7)^2+(G7/F7)^2)*H7"
$ws.Range("I8").Formula = "=SQRT((C8/B8)^2+(E8/D8)^2+(G8/F8)^2)*H8"

# Update selection to I2
$ws.Range("I2").Select()
